# Simulator full-month coverage: populate Rate/Total figures that were
# previously left at 0 now that the timesheet covers the full month.

$wb = $excel.ActiveWorkbook

# --- "Weekly Timesheet" sheet -------------------------------------------------
$ws1 = $wb.Worksheets.Item("Weekly Timesheet")

# Daily rows: Rate (E) / Total (F) for rows 2-5
foreach ($r in 2..5) {
    $ws1.Cells.Item($r, 5).Value = 150
    $ws1.Cells.Item($r, 6).Value = 1200
}

# Subtotal/summary rows: Total (F) only
$ws1.Range("F7").Value = 4800
$ws1.Range("F11").Value = 4800
$ws1.Range("F12").Value = 4800

# --- "Jason Schema" sheet ------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Jason Schema")

# Rate (F) / Total (G) for rows 2-5
foreach ($r in 2..5) {
    $ws2.Cells.Item($r, 6).Value = 150
    $ws2.Cells.Item($r, 7).Value = 1200
}
